# Apply "add auto change mde" edit:
# Insert 9 new data rows at the top of the results table (rows 6-14) on the
# "java" worksheet, pushing the previous rows 6-13 down to rows 15-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("java")

# Insert 9 blank rows before the existing row 6 (shifts old rows 6-13 -> 15-22)
$ws.Range("A6:J14").EntireRow.Insert()

# New rows of data to populate at rows 6-14
$newRows = @(
    @(6,  0, 0,     0,      " 9294.6/s ",  0,                   0,                   "short", "75%Success",  "apm", "10"),
    @(7,  0, 17625, 219236, " 4733.5/s ",  1530.547913558057,  60.89500811910357,  "short", "70%Success",  "apm", "10"),
    @(8,  0, 7860,  437392, " 13840.2/s ", 5048.718359723461,  19.02445403665375,  "short", "100%Success", "apm", "10"),
    @(9,  0, 8415,  509238, " 15535.0/s ", 5661.700839115933,  25.66828477057825,  "short", "95%Success",  "apm", "1"),
    @(10, 0, 25170, 431379, " 12395.2/s ", 4526.671226067946,  31.13852551932311,  "short", "90%Success",  "apm", "1"),
    @(11, 0, 17610, 396499, " 6599.1/s ",  2404.019495873922,  34.71681139170617,  "short", "85%Success",  "apm", "1"),
    @(12, 0, 16695, 356787, " 4138.2/s ",  1506.000210839302,  40.23573728863409,  "short", "80%Success",  "apm", "1"),
    @(13, 0, 14235, 325359, " 5193.6/s ",  1892.590185189449,  45.08529962287786,  "short", "75%Success",  "apm", "1"),
    @(14, 0, 14550, 295284, " 3141.5/s ",  1144.051285127713,  50.49084271413308,  "short", "70%Success",  "apm", "1")
)

foreach ($row in $newRows) {
    $r  = $row[0]
    $a  = $row[1]
    $b  = $row[2]
    $c  = $row[3]
    $d  = $row[4]
    $e  = $row[5]
    $f  = $row[6]
    $g  = $row[7]
    $h  = $row[8]
    $i  = $row[9]
    $j  = $row[10]

    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i

    # Column J values are stored as text in the sheet (e.g. "10", "1"),
    # force text formatting so the numeric-looking string isn't coerced
    # into a real number.
    $ws.Cells.Item($r, 10).NumberFormat = "@"
    $ws.Cells.Item($r, 10).Value = $j
}
